$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 22.95531766666667
$ws.Range("H2").Value = 68.865953
$ws.Range("I2").Value = 0.1720020945576478
$ws.Range("J2").Value = 0.1720020945576478
$ws.Range("M2").Value = 5.854382333333334
$ws.Range("N2").Value = 17.563147
$ws.Range("O2").Value = 0.1730451459016118
$ws.Range("P2").Value = 0.1730451459016118
$ws.Range("Q2").Value = 134.3892062037879
$ws.Range("R2").Value = 1209.502855834091
$ws.Range("S2").Value = 0.029764127548111
$ws.Range("T2").Value = 0.029764127548111

$ws.Range("G3").Value = 22.95531766666667
$ws.Range("H3").Value = 68.865953
$ws.Range("I3").Value = 0.1720020945576478
$ws.Range("J3").Value = 0.1720020945576478
$ws.Range("O3").Value = 0.1208497063316524
$ws.Range("P3").Value = 0.1208497063316525
$ws.Range("Q3").Value = 93.85352024324177
$ws.Range("R3").Value = 844.6816821891761
$ws.Range("S3").Value = 0.02078640261572086
$ws.Range("T3").Value = 0.02078640261572086

$ws.Range("G4").Value = 22.95531766666667
$ws.Range("H4").Value = 68.865953
$ws.Range("I4").Value = 0.1720020945576478
$ws.Range("J4").Value = 0.1720020945576478
$ws.Range("M4").Value = 11.64342866666667
$ws.Range("N4").Value = 34.930286
$ws.Range("O4").Value = 0.3441590756630932
$ws.Range("P4").Value = 0.3441590756630932
$ws.Range("Q4").Value = 267.2786037725064
$ws.Range("R4").Value = 2405.507433952558
$ws.Range("S4").Value = 0.05919608187507603
$ws.Range("T4").Value = 0.05919608187507603

$ws.Range("G5").Value = 22.95531766666667
$ws.Range("H5").Value = 68.865953
$ws.Range("I5").Value = 0.1720020945576478
$ws.Range("J5").Value = 0.1720020945576478
$ws.Range("M5").Value = 0.1645376666666667
$ws.Range("N5").Value = 0.493613
$ws.Range("O5").Value = 0.004863441250245888
$ws.Range("P5").Value = 0.004863441250245888
$ws.Range("Q5").Value = 3.777014406465445
$ws.Range("R5").Value = 33.99312965818901
$ws.Range("S5").Value = 0.0008365220818003582
$ws.Range("T5").Value = 0.0008365220818003582

$ws.Range("G6").Value = 22.95531766666667
$ws.Range("H6").Value = 68.865953
$ws.Range("I6").Value = 0.1720020945576478
$ws.Range("J6").Value = 0.1720020945576478
$ws.Range("M6").Value = 12.08065233333333
$ws.Range("N6").Value = 36.241957
$ws.Range("O6").Value = 0.3570826308533967
$ws.Range("P6").Value = 0.3570826308533967
$ws.Range("Q6").Value = 277.3152119322246
$ws.Range("R6").Value = 2495.836907390021
$ws.Range("S6").Value = 0.0614189604369396
$ws.Range("T6").Value = 0.0614189604369396

$ws.Range("I7").Value = 0.4661646602805707
$ws.Range("J7").Value = 0.4661646602805707
$ws.Range("M7").Value = 5.854382333333334
$ws.Range("N7").Value = 17.563147
$ws.Range("O7").Value = 0.1730451459016118
$ws.Range("P7").Value = 0.1730451459016118
$ws.Range("Q7").Value = 364.2252079341252
$ws.Range("R7").Value = 3278.026871407127
$ws.Range("S7").Value = 0.08066753165242667
$ws.Range("T7").Value = 0.08066753165242667

$ws.Range("I8").Value = 0.4661646602805707
$ws.Range("J8").Value = 0.4661646602805707
$ws.Range("O8").Value = 0.1208497063316524
$ws.Range("P8").Value = 0.1208497063316525
$ws.Range("S8").Value = 0.0563358622971015
$ws.Range("T8").Value = 0.0563358622971015

$ws.Range("I9").Value = 0.4661646602805707
$ws.Range("J9").Value = 0.4661646602805707
$ws.Range("M9").Value = 11.64342866666667
$ws.Range("N9").Value = 34.930286
$ws.Range("O9").Value = 0.3441590756630932
$ws.Range("P9").Value = 0.3441590756630932
$ws.Range("Q9").Value = 724.385594537725
$ws.Range("R9").Value = 6519.470350839525
$ws.Range("S9").Value = 0.1604347985889611
$ws.Range("T9").Value = 0.1604347985889611

$ws.Range("I10").Value = 0.4661646602805707
$ws.Range("J10").Value = 0.4661646602805707
$ws.Range("M10").Value = 0.1645376666666667
$ws.Range("N10").Value = 0.493613
$ws.Range("O10").Value = 0.004863441250245888
$ws.Range("P10").Value = 0.004863441250245888
$ws.Range("Q10").Value = 10.23656509644811
$ws.Range("R10").Value = 92.129085868033
$ws.Range("S10").Value = 0.002267164438215388
$ws.Range("T10").Value = 0.002267164438215388

$ws.Range("I11").Value = 0.4661646602805707
$ws.Range("J11").Value = 0.4661646602805707
$ws.Range("M11").Value = 12.08065233333333
$ws.Range("N11").Value = 36.241957
$ws.Range("O11").Value = 0.3570826308533967
$ws.Range("P11").Value = 0.3570826308533967
$ws.Range("Q11").Value = 751.5870774334818
$ws.Range("R11").Value = 6764.283696901336
$ws.Range("S11").Value = 0.1664593033038661
$ws.Range("T11").Value = 0.1664593033038661

$ws.Range("G12").Value = 7.783044333333334
$ws.Range("H12").Value = 23.349133
$ws.Range("I12").Value = 0.05831763893698088
$ws.Range("J12").Value = 0.05831763893698089
$ws.Range("M12").Value = 5.854382333333334
$ws.Range("N12").Value = 17.563147
$ws.Range("O12").Value = 0.1730451459016118
$ws.Range("P12").Value = 0.1730451459016118
$ws.Range("Q12").Value = 45.56491724461678
$ws.Range("R12").Value = 410.0842552015511
$ws.Range("S12").Value = 0.01009158433848737
$ws.Range("T12").Value = 0.01009158433848738

$ws.Range("G13").Value = 7.783044333333334
$ws.Range("H13").Value = 23.349133
$ws.Range("I13").Value = 0.05831763893698088
$ws.Range("J13").Value = 0.05831763893698089
$ws.Range("O13").Value = 0.1208497063316524
$ws.Range("P13").Value = 0.1208497063316525
$ws.Range("Q13").Value = 31.82121543685955
$ws.Range("R13").Value = 286.390938931736
$ws.Range("S13").Value = 0.007047669539489479
$ws.Range("T13").Value = 0.007047669539489481

$ws.Range("G14").Value = 7.783044333333334
$ws.Range("H14").Value = 23.349133
$ws.Range("I14").Value = 0.05831763893698088
$ws.Range("J14").Value = 0.05831763893698089
$ws.Range("M14").Value = 11.64342866666667
$ws.Range("N14").Value = 34.930286
$ws.Range("O14").Value = 0.3441590756630932
$ws.Range("P14").Value = 0.3441590756630932
$ws.Range("Q14").Value = 90.62132150467087
$ws.Range("R14").Value = 815.591893542038
$ws.Range("S14").Value = 0.02007054471140535
$ws.Range("T14").Value = 0.02007054471140535

$ws.Range("G15").Value = 7.783044333333334
$ws.Range("H15").Value = 23.349133
$ws.Range("I15").Value = 0.05831763893698088
$ws.Range("J15").Value = 0.05831763893698089
$ws.Range("M15").Value = 0.1645376666666667
$ws.Range("N15").Value = 0.493613
$ws.Range("O15").Value = 0.004863441250245888
$ws.Range("P15").Value = 0.004863441250245888
$ws.Range("Q15").Value = 1.280603954169889
$ws.Range("R15").Value = 11.525435587529
$ws.Range("S15").Value = 0.0002836244108230586
$ws.Range("T15").Value = 0.0002836244108230586

$ws.Range("G16").Value = 7.783044333333334
$ws.Range("H16").Value = 23.349133
$ws.Range("I16").Value = 0.05831763893698088
$ws.Range("J16").Value = 0.05831763893698089
$ws.Range("M16").Value = 12.08065233333333
$ws.Range("N16").Value = 36.241957
$ws.Range("O16").Value = 0.3570826308533967
$ws.Range("P16").Value = 0.3570826308533967
$ws.Range("Q16").Value = 94.02425268592012
$ws.Range("R16").Value = 846.2182741732811
$ws.Range("S16").Value = 0.02082421593677562
$ws.Range("T16").Value = 0.02082421593677562

$ws.Range("G17").Value = 30.44016466666666
$ws.Range("H17").Value = 91.320494
$ws.Range("I17").Value = 0.2280853681650076
$ws.Range("J17").Value = 0.2280853681650076
$ws.Range("M17").Value = 5.854382333333334
$ws.Range("N17").Value = 17.563147
$ws.Range("O17").Value = 0.1730451459016118
$ws.Range("P17").Value = 0.1730451459016118
$ws.Range("Q17").Value = 178.2083622482909
$ws.Range("R17").Value = 1603.875260234618
$ws.Range("S17").Value = 0.03946906581213659
$ws.Range("T17").Value = 0.03946906581213659

$ws.Range("G18").Value = 30.44016466666666
$ws.Range("H18").Value = 91.320494
$ws.Range("I18").Value = 0.2280853681650076
$ws.Range("J18").Value = 0.2280853681650076
$ws.Range("O18").Value = 0.1208497063316524
$ws.Range("P18").Value = 0.1208497063316525
$ws.Range("Q18").Value = 124.4555467380498
$ws.Range("R18").Value = 1120.099920642448
$ws.Range("S18").Value = 0.027564049761288
$ws.Range("T18").Value = 0.027564049761288

$ws.Range("G19").Value = 30.44016466666666
$ws.Range("H19").Value = 91.320494
$ws.Range("I19").Value = 0.2280853681650076
$ws.Range("J19").Value = 0.2280853681650076
$ws.Range("M19").Value = 11.64342866666667
$ws.Range("N19").Value = 34.930286
$ws.Range("O19").Value = 0.3441590756630932
$ws.Range("P19").Value = 0.3441590756630932
$ws.Range("Q19").Value = 354.4278858979204
$ws.Range("R19").Value = 3189.850973081283
$ws.Range("S19").Value = 0.07849764947994531
$ws.Range("T19").Value = 0.07849764947994531

$ws.Range("G20").Value = 30.44016466666666
$ws.Range("H20").Value = 91.320494
$ws.Range("I20").Value = 0.2280853681650076
$ws.Range("J20").Value = 0.2280853681650076
$ws.Range("M20").Value = 0.1645376666666667
$ws.Range("N20").Value = 0.493613
$ws.Range("O20").Value = 0.004863441250245888
$ws.Range("P20").Value = 0.004863441250245888
$ws.Range("Q20").Value = 5.008553667202444
$ws.Range("R20").Value = 45.076983004822
$ws.Range("S20").Value = 0.001109279788111218
$ws.Range("T20").Value = 0.001109279788111218

$ws.Range("G21").Value = 30.44016466666666
$ws.Range("H21").Value = 91.320494
$ws.Range("I21").Value = 0.2280853681650076
$ws.Range("J21").Value = 0.2280853681650076
$ws.Range("M21").Value = 12.08065233333333
$ws.Range("N21").Value = 36.241957
$ws.Range("O21").Value = 0.3570826308533967
$ws.Range("P21").Value = 0.3570826308533967
$ws.Range("Q21").Value = 367.7370463074175
$ws.Range("R21").Value = 3309.633416766758
$ws.Range("S21").Value = 0.08144532332352648
$ws.Range("T21").Value = 0.08144532332352648

$ws.Range("G22").Value = 10.06688366666667
$ws.Range("H22").Value = 30.200651
$ws.Range("I22").Value = 0.07543023805979308
$ws.Range("J22").Value = 0.07543023805979308
$ws.Range("M22").Value = 5.854382333333334
$ws.Range("N22").Value = 17.563147
$ws.Range("O22").Value = 0.1730451459016118
$ws.Range("P22").Value = 0.1730451459016118
$ws.Range("Q22").Value = 58.93538588985523
$ws.Range("R22").Value = 530.418473008697
$ws.Range("S22").Value = 0.01305283655045021
$ws.Range("T22").Value = 0.01305283655045021

$ws.Range("G23").Value = 10.06688366666667
$ws.Range("H23").Value = 30.200651
$ws.Range("I23").Value = 0.07543023805979308
$ws.Range("J23").Value = 0.07543023805979308
$ws.Range("O23").Value = 0.1208497063316524
$ws.Range("P23").Value = 0.1208497063316525
$ws.Range("Q23").Value = 41.15876258893245
$ws.Range("R23").Value = 370.428863300392
$ws.Range("S23").Value = 0.009115722118052628
$ws.Range("T23").Value = 0.009115722118052628

$ws.Range("G24").Value = 10.06688366666667
$ws.Range("H24").Value = 30.200651
$ws.Range("I24").Value = 0.07543023805979308
$ws.Range("J24").Value = 0.07543023805979308
$ws.Range("M24").Value = 11.64342866666667
$ws.Range("N24").Value = 34.930286
$ws.Range("O24").Value = 0.3441590756630932
$ws.Range("P24").Value = 0.3441590756630932
$ws.Range("Q24").Value = 117.2130418684651
$ws.Range("R24").Value = 1054.917376816186
$ws.Range("S24").Value = 0.02596000100770546
$ws.Range("T24").Value = 0.02596000100770546

$ws.Range("G25").Value = 10.06688366666667
$ws.Range("H25").Value = 30.200651
$ws.Range("I25").Value = 0.07543023805979308
$ws.Range("J25").Value = 0.07543023805979308
$ws.Range("M25").Value = 0.1645376666666667
$ws.Range("N25").Value = 0.493613
$ws.Range("O25").Value = 0.004863441250245888
$ws.Range("P25").Value = 0.004863441250245888
$ws.Range("Q25").Value = 1.656381549118111
$ws.Range("R25").Value = 14.907433942063
$ws.Range("S25").Value = 0.000366850531295865
$ws.Range("T25").Value = 0.000366850531295865

$ws.Range("G26").Value = 10.06688366666667
$ws.Range("H26").Value = 30.200651
$ws.Range("I26").Value = 0.07543023805979308
$ws.Range("J26").Value = 0.07543023805979308
$ws.Range("M26").Value = 12.08065233333333
$ws.Range("N26").Value = 36.241957
$ws.Range("O26").Value = 0.3570826308533967
$ws.Range("P26").Value = 0.3570826308533967
$ws.Range("Q26").Value = 121.6145216571119
$ws.Range("R26").Value = 1094.530694914007
$ws.Range("S26").Value = 0.02693482785228893
$ws.Range("T26").Value = 0.02693482785228893
